# Weekly price update: insert a new daily record for
# "Feria Lagunitas de Puerto Montt - Zanahoria" as row 363, shifting the
# existing rows 363:388 down to 364:389.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 363:388 down one row, freeing up row 363 for the new record.
$ws.Rows.Item(363).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(363, 1).Value = 4
$ws.Cells.Item(363, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(363, 3).Value = "Los Lagos"
$ws.Cells.Item(363, 4).Value = 44746
$ws.Cells.Item(363, 5).Value = 10
$ws.Cells.Item(363, 6).Value = 100114013
$ws.Cells.Item(363, 7).Value = "Zanahoria"
$ws.Cells.Item(363, 8).Value = "Sin especificar"
$ws.Cells.Item(363, 9).Value = "Primera"
$ws.Cells.Item(363, 10).Value = 250
$ws.Cells.Item(363, 11).Value = 10000
$ws.Cells.Item(363, 12).Value = 10000
$ws.Cells.Item(363, 13).Value = 10000
$ws.Cells.Item(363, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(363, 15).Value = "Chillán"
$ws.Cells.Item(363, 16).Value = 500
$ws.Cells.Item(363, 17).Value = 20
$ws.Cells.Item(363, 18).Value = "Hortaliza"
